# Refactored and QoL changes
# Reset the sample/demo row into a reusable template:
#  - Row 1 becomes column headers (text labels) instead of sample data.
#  - Row 2 gets a single placeholder cell telling the user where to start
#    entering real data.
#  - Selection moves to A2 (the first data-entry cell) and the window
#    scroll position is nudged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Roll no.", "CGPA", "Subject Code", "Grade", "SGPA 1", "SGPA 2", "SGPA 3", "SGPA 4", "SGPA 5", "SGPA 6", "SGPA 7", "SGPA 8")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("A2").Value = "<Start inserting data from here>"

[void]$ws.Range("A2").Select()

# Best-effort: nudge the remembered window scroll position (headless runs
# have no real screen, so this may be a no-op, but mirrors the authored
# workbookView xWindow/yWindow change).
$win = $excel.ActiveWindow
$win.Left = 9012
$win.Top = 4140
